# Update industrial biomass extension.
# - Set the selection on the "petroleum" sheet to the header row A1:C1.
# - Add a new "other_biomass" worksheet at the end of the workbook, populate
#   it with the sector / X1850 data, size column A, and make it the active
#   (selected) sheet/tab with B1 selected.
# - Give the "petroleum" sheet an explicit (portrait) page setup, matching
#   the "coal" sheet.

$wb = $excel.ActiveWorkbook

# --- petroleum sheet: update selection + page setup -----------------------
$petroleum = $wb.Worksheets.Item("petroleum")
$petroleum.Range("A1:C1").Select() | Out-Null
$petroleum.PageSetup.Orientation = 1

# --- new sheet: other_biomass ---------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$biomass = $wb.Worksheets.Add($null, $lastSheet)
$biomass.Name = "other_biomass"

# Column A width to match the source sheets' formatting.
$biomass.Columns.Item(1).ColumnWidth = 23.8

# Page setup: portrait orientation with the same margins used by the other
# sheets in the workbook (1in top/bottom, 0.75in left/right, 0.5in header
# and footer).
$biomassPageSetup = $biomass.PageSetup
$biomassPageSetup.Orientation = 1
$biomassPageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$biomassPageSetup.RightMargin = $excel.InchesToPoints(0.75)
$biomassPageSetup.TopMargin = $excel.InchesToPoints(1)
$biomassPageSetup.BottomMargin = $excel.InchesToPoints(1)
$biomassPageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$biomassPageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# Row 9 (new sector string) must be written before B1 (new "X1850" header
# string) so the shared-string table gets the same ordering as the source
# workbook (navigation sector string allocated right before the X1850
# header string).
$biomass.Range("A1").Value = "sector"
$biomass.Range("A2").Value = "1A1a_Electricity-autoproducer"
$biomass.Range("A3").Value = "1A1a_Electricity-public"
$biomass.Range("A4").Value = "1A1a_Heat-production"
$biomass.Range("A5").Value = "1A3ai_International-aviation"
$biomass.Range("A6").Value = "1A3aii_Domestic-aviation"
$biomass.Range("A7").Value = "1A3b_Road"
$biomass.Range("A8").Value = "1A3c_Rail"
$biomass.Range("A9").Value = "1A3dii_Domestic-navigation"
$biomass.Range("A10").Value = "1A3eii_Other-transp"
$biomass.Range("A11").Value = "1A4a_Commercial-institutional"
$biomass.Range("A12").Value = "1A4c_Agriculture-forestry-fishing"
$biomass.Range("A13").Value = "1A5_Other-unspecified"

$biomass.Range("B1").Value = "X1850"
$biomass.Range("B2").Value = 0
$biomass.Range("B3").Value = 0
$biomass.Range("B4").Value = 0
$biomass.Range("B5").Value = 0
$biomass.Range("B6").Value = 0
$biomass.Range("B7").Value = 0
$biomass.Range("B8").Value = 0
$biomass.Range("B9").Value = 0
$biomass.Range("B10").Value = 0
$biomass.Range("B11").Value = 0
$biomass.Range("B12").Value = 0
$biomass.Range("B13").Value = 1

# Make the new sheet the active tab/selection, matching the authored file.
$biomass.Range("B1").Select() | Out-Null
